$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values regenerated for sval data (filtered save games)
$data = @{
    2 = @(0.1169995834814548, 0.04103571897497393, 0.7210945179870265, 13.86384647080068, 0)
    3 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0)
    4 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1)
    5 = @(0.1169995834814548, 0.04103571897497393, 189.6080260415259, 13.86384647080068, 1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]
    $f = $vals[4]
    $g = $b + $c + $d + $e

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
}
